# Atualizacao de bases das ligas, do dia: 28-05-2024 as 19:13
# This script re-syncs a handful of match rows in the "Romania Liga I"
# sheet so that each row's data (id, HomeTeam, AwayTeam, scores, odds, ...)
# reflects the latest scrape, while leaving the row's running index (col A)
# and match Date (col D) untouched.
#
# Concretely, the data (columns B,C,E:AD) held by some rows got reshuffled
# among themselves (the underlying matches kept their row position/date,
# but the fixture info moved to a different row than before).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Block 1: rows 236-239 -------------------------------------------------
# Capture the "before" payload of each row first (so source values are not
# clobbered before they've been copied elsewhere).
$row236 = $ws.Range("B236:AD236").Value2
$row237 = $ws.Range("B237:AD237").Value2
$row238 = $ws.Range("B238:AD238").Value2
$row239 = $ws.Range("B239:AD239").Value2

# New row236 gets what used to be row239's data, etc.
$ws.Range("B236:AD236").Value2 = $row239
$ws.Range("B237:AD237").Value2 = $row236
$ws.Range("B238:AD238").Value2 = $row237
$ws.Range("B239:AD239").Value2 = $row238

# ---- Block 2: rows 309-316 -------------------------------------------------
$row309 = $ws.Range("B309:AD309").Value2
$row310 = $ws.Range("B310:AD310").Value2
$row311 = $ws.Range("B311:AD311").Value2
$row312 = $ws.Range("B312:AD312").Value2
$row313 = $ws.Range("B313:AD313").Value2
$row315 = $ws.Range("B315:AD315").Value2
$row316 = $ws.Range("B316:AD316").Value2
# row314 is unchanged, no need to capture/restore it

$ws.Range("B309:AD309").Value2 = $row311
$ws.Range("B310:AD310").Value2 = $row309
$ws.Range("B311:AD311").Value2 = $row312
$ws.Range("B312:AD312").Value2 = $row313
$ws.Range("B313:AD313").Value2 = $row310
$ws.Range("B315:AD315").Value2 = $row316
$ws.Range("B316:AD316").Value2 = $row315
